$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (subject ids) for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update "CON" row (row 2) values for columns B:E
$ws.Range("B2").Value = [double]"2.9867720791489222E-2"
$ws.Range("C2").Value = [double]"-0.1223857107086631"
$ws.Range("D2").Value = [double]"9.5783367102774079E-3"
$ws.Range("E2").Value = [double]"2.4770903962902775E-2"

# Update "STR" row (row 3) values for columns B:E
$ws.Range("B3").Value = [double]"4.5045087332597329E-2"
$ws.Range("C3").Value = [double]"0.11486536248197167"
$ws.Range("D3").Value = [double]"3.0193810726400653E-2"
$ws.Range("E3").Value = [double]"4.8525950261542029E-2"

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
